$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.021598327441391
$ws.Cells.Item(2, 4).Value = 1.032306816338145
$ws.Cells.Item(2, 5).Value = 1.02245750295812
$ws.Cells.Item(2, 6).Value = 1.040884514804083
$ws.Cells.Item(2, 9).Value = 1.031202846778234
$ws.Cells.Item(2, 10).Value = 1.026788892010296
$ws.Cells.Item(2, 11).Value = 1.035112425782007
$ws.Cells.Item(2, 12).Value = 1.025291826926038
$ws.Cells.Item(2, 13).Value = 1.043665597723298
$ws.Cells.Item(2, 14).Value = 1.01300092151919

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.022487249434818
$ws.Cells.Item(3, 4).Value = 1.03302777819261
$ws.Cells.Item(3, 5).Value = 1.023209632348831
$ws.Cells.Item(3, 6).Value = 1.041840884720024
$ws.Cells.Item(3, 9).Value = 1.031367730561927
$ws.Cells.Item(3, 10).Value = 1.027315817055318
$ws.Cells.Item(3, 11).Value = 1.035642367475331
$ws.Cells.Item(3, 12).Value = 1.025850775610291
$ws.Cells.Item(3, 13).Value = 1.044432098885888
$ws.Cells.Item(3, 14).Value = 1.013175609733828

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.023062751158941
$ws.Cells.Item(4, 4).Value = 1.0334940837698
$ws.Cells.Item(4, 5).Value = 1.023696983382513
$ws.Cells.Item(4, 6).Value = 1.042460030013461
$ws.Cells.Item(4, 9).Value = 1.031472553555908
$ws.Cells.Item(4, 10).Value = 1.027656469080816
$ws.Cells.Item(4, 11).Value = 1.035984402235043
$ws.Cells.Item(4, 12).Value = 1.026212464552057
$ws.Cells.Item(4, 13).Value = 1.044927766986196
$ws.Cells.Item(4, 14).Value = 1.013288517074623

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.023304765068553
$ws.Cells.Item(5, 4).Value = 1.033690067995438
$ws.Cells.Item(5, 5).Value = 1.023902025385338
$ws.Cells.Item(5, 6).Value = 1.0427203914116
$ws.Cells.Item(5, 9).Value = 1.031516173172
$ws.Cells.Item(5, 10).Value = 1.027799605427286
$ws.Cells.Item(5, 11).Value = 1.03612798340788
$ws.Cells.Item(5, 12).Value = 1.026364520265464
$ws.Cells.Item(5, 13).Value = 1.045136070385026
$ws.Cells.Item(5, 14).Value = 1.013335952388048

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.023345404514697
$ws.Cells.Item(6, 4).Value = 1.033722971593913
$ws.Cells.Item(6, 5).Value = 1.023936462156957
$ws.Cells.Item(6, 6).Value = 1.042764111456446
$ws.Cells.Item(6, 9).Value = 1.031523470804728
$ws.Cells.Item(6, 10).Value = 1.027823634293306
$ws.Cells.Item(6, 11).Value = 1.036152078953427
$ws.Cells.Item(6, 12).Value = 1.026390051154402
$ws.Cells.Item(6, 13).Value = 1.045171040975883
$ws.Cells.Item(6, 14).Value = 1.013343915160878

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.023065984675335
$ws.Cells.Item(7, 4).Value = 1.033496702721584
$ws.Cells.Item(7, 5).Value = 1.023699722538448
$ws.Cells.Item(7, 6).Value = 1.042463508691373
$ws.Cells.Item(7, 9).Value = 1.031473138164594
$ws.Cells.Item(7, 10).Value = 1.027658381966863
$ws.Cells.Item(7, 11).Value = 1.035986321601391
$ws.Cells.Item(7, 12).Value = 1.026214496322765
$ws.Cells.Item(7, 13).Value = 1.044930550644722
$ws.Cells.Item(7, 14).Value = 1.013289151029629

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.021898677950853
$ws.Cells.Item(8, 4).Value = 1.032550510265555
$ws.Cells.Item(8, 5).Value = 1.022711548355516
$ws.Cells.Item(8, 6).Value = 1.041207659405782
$ws.Cells.Item(8, 9).Value = 1.031258956290477
$ws.Cells.Item(8, 10).Value = 1.026967031048733
$ws.Cells.Item(8, 11).Value = 1.035291702057167
$ws.Cells.Item(8, 12).Value = 1.025480723133316
$ws.Cells.Item(8, 13).Value = 1.043924703605722
$ws.Cells.Item(8, 14).Value = 1.013059984413955

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.019844158180194
$ws.Cells.Item(9, 4).Value = 1.030881697635309
$ws.Cells.Item(9, 5).Value = 1.020975479554507
$ws.Cells.Item(9, 6).Value = 1.038997132615315
$ws.Cells.Item(9, 9).Value = 1.030867270316939
$ws.Cells.Item(9, 10).Value = 1.025746502674897
$ws.Cells.Item(9, 11).Value = 1.034061062874545
$ws.Cells.Item(9, 12).Value = 1.024187861032249
$ws.Cells.Item(9, 13).Value = 1.042149964259018
$ws.Cells.Item(9, 14).Value = 1.01265520293103

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.018476175989034
$ws.Cells.Item(10, 4).Value = 1.029768238951514
$ws.Cells.Item(10, 5).Value = 1.019821698244444
$ws.Cells.Item(10, 6).Value = 1.037525175894082
$ws.Cells.Item(10, 9).Value = 1.030596598697321
$ws.Cells.Item(10, 10).Value = 1.024931345978082
$ws.Cells.Item(10, 11).Value = 1.033236251636873
$ws.Cells.Item(10, 12).Value = 1.023326113629685
$ws.Cells.Item(10, 13).Value = 1.040965334030676
$ws.Cells.Item(10, 14).Value = 1.012384725057726

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.017884241160452
$ws.Cells.Item(11, 4).Value = 1.029285900632457
$ws.Cells.Item(11, 5).Value = 1.019322969923574
$ws.Cells.Item(11, 6).Value = 1.036888230111256
$ws.Cells.Item(11, 9).Value = 1.030477139945506
$ws.Cells.Item(11, 10).Value = 1.024578037931181
$ws.Cells.Item(11, 11).Value = 1.032878074245727
$ws.Cells.Item(11, 12).Value = 1.022953019402633
$ws.Cells.Item(11, 13).Value = 1.040452043970046
$ws.Cells.Item(11, 14).Value = 1.012267461876986

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.017664432949183
$ws.Cells.Item(12, 4).Value = 1.029106709650961
$ws.Cells.Item(12, 5).Value = 1.019137851679299
$ws.Cells.Item(12, 6).Value = 1.036651704758743
$ws.Cells.Item(12, 9).Value = 1.030432429414136
$ws.Cells.Item(12, 10).Value = 1.024446753530208
$ws.Cells.Item(12, 11).Value = 1.032744877836535
$ws.Cells.Item(12, 12).Value = 1.022814443906262
$ws.Cells.Item(12, 13).Value = 1.04026133580071
$ws.Cells.Item(12, 14).Value = 1.012223883766336

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.017711579700388
$ws.Cells.Item(13, 4).Value = 1.02914514801552
$ws.Cells.Item(13, 5).Value = 1.019177554189765
$ws.Cells.Item(13, 6).Value = 1.036702437300952
$ws.Cells.Item(13, 9).Value = 1.030442035268937
$ws.Cells.Item(13, 10).Value = 1.024474916729356
$ws.Cells.Item(13, 11).Value = 1.032773455840427
$ws.Cells.Item(13, 12).Value = 1.022844168427105
$ws.Cells.Item(13, 13).Value = 1.040302245564483
$ws.Cells.Item(13, 14).Value = 1.012233232379074

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.017866070459082
$ws.Cells.Item(14, 4).Value = 1.029271089234071
$ws.Cells.Item(14, 5).Value = 1.019307665297487
$ws.Cells.Item(14, 6).Value = 1.036868677524526
$ws.Cells.Item(14, 9).Value = 1.030473451054663
$ws.Cells.Item(14, 10).Value = 1.024567186937151
$ws.Cells.Item(14, 11).Value = 1.03286706731112
$ws.Cells.Item(14, 12).Value = 1.022941564533046
$ws.Cells.Item(14, 13).Value = 1.040436280973266
$ws.Cells.Item(14, 14).Value = 1.012263860130092

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.017961265659971
$ws.Cells.Item(15, 4).Value = 1.029348681970222
$ws.Cells.Item(15, 5).Value = 1.019387848535065
$ws.Cells.Item(15, 6).Value = 1.036971112224014
$ws.Cells.Item(15, 9).Value = 1.030492762559817
$ws.Cells.Item(15, 10).Value = 1.024624031053489
$ws.Cells.Item(15, 11).Value = 1.032924724138553
$ws.Cells.Item(15, 12).Value = 1.023001574622215
$ws.Cells.Item(15, 13).Value = 1.04051885810685
$ws.Cells.Item(15, 14).Value = 1.012282728080825

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.018515469493865
$ws.Cells.Item(16, 4).Value = 1.029800246000744
$ws.Cells.Item(16, 5).Value = 1.01985481558385
$ws.Cells.Item(16, 6).Value = 1.0375674568809
$ws.Cells.Item(16, 9).Value = 1.030604479309062
$ws.Cells.Item(16, 10).Value = 1.024954786794792
$ws.Cells.Item(16, 11).Value = 1.033260001117516
$ws.Cells.Item(16, 12).Value = 1.023350875770777
$ws.Cells.Item(16, 13).Value = 1.040999392404692
$ws.Cells.Item(16, 14).Value = 1.012392504415792

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.018863217571
$ws.Cells.Item(17, 4).Value = 1.030083446801997
$ws.Cells.Item(17, 5).Value = 1.020147964724395
$ws.Cells.Item(17, 6).Value = 1.037941641878063
$ws.Cells.Item(17, 9).Value = 1.030673952754408
$ws.Cells.Item(17, 10).Value = 1.025162170855824
$ws.Cells.Item(17, 11).Value = 1.033470036928737
$ws.Cells.Item(17, 12).Value = 1.02356999679545
$ws.Cells.Item(17, 13).Value = 1.041300729472308
$ws.Cells.Item(17, 14).Value = 1.012461325785828

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.019066092675431
$ws.Cells.Item(18, 4).Value = 1.030248613316472
$ws.Cells.Item(18, 5).Value = 1.020319037174903
$ws.Cells.Item(18, 6).Value = 1.038159938201945
$ws.Cells.Item(18, 9).Value = 1.030714257704842
$ws.Cells.Item(18, 10).Value = 1.025283101616115
$ws.Cells.Item(18, 11).Value = 1.033592447932883
$ws.Cells.Item(18, 12).Value = 1.023697810948708
$ws.Cells.Item(18, 13).Value = 1.041476461649148
$ws.Cells.Item(18, 14).Value = 1.012501454162513

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.019135274516479
$ws.Cells.Item(19, 4).Value = 1.030304927438792
$ws.Cells.Item(19, 5).Value = 1.02037738258088
$ws.Cells.Item(19, 6).Value = 1.038234378441942
$ws.Cells.Item(19, 9).Value = 1.03072796369749
$ws.Cells.Item(19, 10).Value = 1.025324330277409
$ws.Cells.Item(19, 11).Value = 1.033634170037974
$ws.Cells.Item(19, 12).Value = 1.023741393013424
$ws.Cells.Item(19, 13).Value = 1.041536376204866
$ws.Cells.Item(19, 14).Value = 1.012515134522372

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.018825903396374
$ws.Cells.Item(20, 4).Value = 1.030053064065223
$ws.Cells.Item(20, 5).Value = 1.020116503946851
$ws.Cells.Item(20, 6).Value = 1.037901491165167
$ws.Cells.Item(20, 9).Value = 1.030666521427244
$ws.Cells.Item(20, 10).Value = 1.025139923899763
$ws.Cells.Item(20, 11).Value = 1.033447512339732
$ws.Cells.Item(20, 12).Value = 1.023546486702995
$ws.Cells.Item(20, 13).Value = 1.041268402223098
$ws.Cells.Item(20, 14).Value = 1.012453943345889

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.017820575040717
$ws.Cells.Item(21, 4).Value = 1.029234003474232
$ws.Cells.Item(21, 5).Value = 1.019269347173012
$ws.Cells.Item(21, 6).Value = 1.036819722125713
$ws.Cells.Item(21, 9).Value = 1.030464209219808
$ws.Cells.Item(21, 10).Value = 1.024540017034716
$ws.Cells.Item(21, 11).Value = 1.032839505290834
$ws.Cells.Item(21, 12).Value = 1.022912883569227
$ws.Cells.Item(21, 13).Value = 1.040396812239045
$ws.Cells.Item(21, 14).Value = 1.012254841607099

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.017188848725022
$ws.Cells.Item(22, 4).Value = 1.028718860213365
$ws.Cells.Item(22, 5).Value = 1.018737468550462
$ws.Cells.Item(22, 6).Value = 1.036139945848893
$ws.Cells.Item(22, 9).Value = 1.030335050956752
$ws.Cells.Item(22, 10).Value = 1.024162542627776
$ws.Cells.Item(22, 11).Value = 1.03245634030274
$ws.Cells.Item(22, 12).Value = 1.022514560755104
$ws.Cells.Item(22, 13).Value = 1.039848523925559
$ws.Cells.Item(22, 14).Value = 1.012129535128863

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.017523703983108
$ws.Cells.Item(23, 4).Value = 1.028991962767866
$ws.Cells.Item(23, 5).Value = 1.019019354711663
$ws.Cells.Item(23, 6).Value = 1.036500272119675
$ws.Cells.Item(23, 9).Value = 1.030403705411558
$ws.Cells.Item(23, 10).Value = 1.024362676030658
$ws.Cells.Item(23, 11).Value = 1.032659546984436
$ws.Cells.Item(23, 12).Value = 1.022725714244361
$ws.Cells.Item(23, 13).Value = 1.040139208526707
$ws.Cells.Item(23, 14).Value = 1.012195974032031

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.018842763943632
$ws.Cells.Item(24, 4).Value = 1.030066792777237
$ws.Cells.Item(24, 5).Value = 1.020130719460112
$ws.Cells.Item(24, 6).Value = 1.037919633418881
$ws.Cells.Item(24, 9).Value = 1.03066987999693
$ws.Cells.Item(24, 10).Value = 1.025149976443195
$ws.Cells.Item(24, 11).Value = 1.033457690538955
$ws.Cells.Item(24, 12).Value = 1.023557109887787
$ws.Cells.Item(24, 13).Value = 1.041283009616445
$ws.Cells.Item(24, 14).Value = 1.012457279195904

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.020375006538528
$ws.Cells.Item(25, 4).Value = 1.031313294023234
$ws.Cells.Item(25, 5).Value = 1.021423667224527
$ws.Cells.Item(25, 6).Value = 1.039568308603529
$ws.Cells.Item(25, 9).Value = 1.030970216772738
$ws.Cells.Item(25, 10).Value = 1.026062302306869
$ws.Cells.Item(25, 11).Value = 1.034379990831764
$ws.Cells.Item(25, 12).Value = 1.024522073411926
$ws.Cells.Item(25, 13).Value = 1.042609042067224
$ws.Cells.Item(25, 14).Value = 1.012759960180472
